# Receding horizon attempt: the heat_pump1 column interactions for
# pv1/bat1/CHP1/pvt1 (row 7) are no longer string-referenced labels;
# replace them with plain numeric zeros like the rest of the row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
